# Bitacora.xlsx — "creacion de modelo join" commit
#
# The log table had an entry in rows 23-24 ("Modelos", dated far in the
# future by mistake, commenting on both the join model and the
# expenses/categorias models at once). The author:
#   1) fixed that entry's date and trimmed its comment down to just the
#      expenses/categorias models, and
#   2) added a new log entry in the next slot (rows 25-26) for the join
#      model work, dated the following day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 25-26 are an already-merged, still-blank entry slot. Copy the
# formatting (fonts/borders/date number format) from the existing entry
# above (rows 23-24) so the new row matches the table's look exactly.
$ws.Range("C23:E24").Copy()
$ws.Range("C25").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# New log entry: "Modelos" / 2022-06-28 / join model comment.
$ws.Range("C25").Value2 = "Modelos"
$ws.Range("D25").Value2 = 44740
$ws.Range("E25").Value2 = "Creacion de modelo join para la consulta en BD"

# Correct the previous entry's date (was a stray far-future typo) and
# narrow its comment to the expenses/categorias models only, now that
# the join-model work has its own row.
$ws.Range("D23").Value2 = 44739
$ws.Range("E23").Value2 = "Creacion de modelos expenses y categorias"

# Match the author's on-save view state: scrolled down a bit further and
# with the new row's neighbourhood selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H25").Select()
